$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing table down by one row to make room for the new
# grouped header row (old row 2 -> row 3, ... old row 8 -> row 9,
# old row 10 -> row 11).
$ws.Rows("2:2").Insert()

# Remove the old standalone "Reell verdi" label row (now row 11) - its
# text is being replaced by the new "Video"/"Radar [m/s]" group headers.
$ws.Rows("11:11").Delete()

# New "Radar [m/s]" / "Video" group header row, centered across the
# merged ranges. Set the text in an order that makes "Video" land
# before "Radar [m/s]" in the shared-string table.
$ws.Range("E2").Value = "Video"
$ws.Range("B2").Value = "Radar [m/s]"

$ws.Range("B2:D2").HorizontalAlignment = -4108
$ws.Range("B2:D2").MergeCells = $true
$ws.Range("E2:G2").HorizontalAlignment = -4108
$ws.Range("E2:G2").MergeCells = $true

# Column headers for the new Video sub-table (row 3), re-using the
# same labels as the existing Radar sub-table.
$ws.Range("E3").Value = "Sakte"
$ws.Range("F3").Value = "kjapp"
$ws.Range("G3").Value = "rygging"

# Sample video measurements.
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

# Average / STD rows for the new Video columns.
$ws.Range("E8").Formula = "=AVERAGE(E4:E7)"
$ws.Range("F8:G8").Formula = "=AVERAGE(F4:F7)"

$ws.Range("E9").Formula = "=STDEV(E4:E7)"
$ws.Range("F9:G9").Formula = "=STDEV(F4:F7)"

# The D column average/std formulas are no longer part of the shared
# formula group after the table was split in two.
$ws.Range("D8").Formula = "=AVERAGE(D4:D7)"
$ws.Range("D9").Formula = "=STDEV(D4:D7)"

$ws.Range("B2:D2").Select()
